$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.4
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("AG2").Value = 8
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 41
$ws.Range("AK2").Value = 34

# Row 3 updates
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
